$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.269.18'
$ws.Range("E2").Value = '  +5.05%  '
$ws.Range("D3").Value = '2.720.52'
$ws.Range("E3").Value = '  +4.52%  '
$origStyle = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = $origStyle
$ws.Range("E4").Value = '  +0.26%  '
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '589.95'
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = '  +1.17%  '
$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '151.05'
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = '  +5.75%  '
$origStyle = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.997'
$ws.Range("D7").Style = $origStyle
$ws.Range("E7").Value = '  -0.18%  '
$ws.Range("E8").Value = '  +1.59%  '
$ws.Range("D9").Value = '2.752.46'
$ws.Range("E9").Value = '  +5.49%  '
$ws.Range("E10").Value = '  +3.45%  '
$ws.Range("E11").Value = '  +7.56%  '
$ws.Range("E12").Value = '  +5.08%  '
$ws.Range("E13").Value = '  +1.49%  '
$ws.Range("D14").Value = '3.210.22'
$ws.Range("E14").Value = '  +4.81%  '
$origStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '26.80'
$ws.Range("D15").Style = $origStyle
$ws.Range("E15").Value = '  +9.01%  '
$ws.Range("D16").Value = '63.180.26'
$ws.Range("E16").Value = '  +4.89%  '
$ws.Range("E17").Value = '  +7.82%  '
$ws.Range("D18").Value = '2.746.59'
$ws.Range("E18").Value = '  +5.39%  '
$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.01'
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = '  +6.07%  '
$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.87'
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = '  +5.61%  '
$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '364.17'
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = '  +4.80%  '
$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.04'
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = '  +2.15%  '
$ws.Range("E23").Value = '  -0.04%  '
$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.536'
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = '  +0.65%  '
$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '65.67'
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = '  +3.05%  '
$ws.Range("E26").Value = '  +4.36%  '
$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.57'
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = '  +7.36%  '
$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.995'
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = '  -0.60%  '
$ws.Range("E29").Value = '  +8.78%  '
$ws.Range("E30").Value = '  +6.87%  '
$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.13'
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = '  +11.38%  '
$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '170.80'
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = '  +1.04%  '
$ws.Range("B33").Value = 'Fetch.AI'
$ws.Range("C33").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.20'
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = '  +20.67%  '
$ws.Range("B34").Value = 'USDe'
$ws.Range("C34").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$origStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.997'
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = '  -0.11%  '
$origStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '20.58'
$ws.Range("D35").Style = $origStyle
$ws.Range("E35").Value = '  +5.98%  '
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$origStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.45'
$ws.Range("D36").Style = $origStyle
$ws.Range("E36").Value = '  +9.18%  '
$ws.Range("B37").Value = 'NEARProtocol'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$origStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.76'
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = '  +12.08%  '
$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.81'
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = '  +9.83%  '
$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.01'
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = '  +19.99%  '
$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '350.90'
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = '  +11.17%  '
$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.28'
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = '  +10.38%  '
$origStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '39.28'
$ws.Range("D42").Style = $origStyle
$ws.Range("E42").Value = '  +2.84%  '
$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '22.46'
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = '  +12.74%  '
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.69'
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = '  +14.55%  '
$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '21.95'
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = '  +10.06%  '
$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0593'
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = '  +8.26%  '
$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '139.62'
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = '  +3.01%  '
$ws.Range("B48").Value = 'Mantle'
$ws.Range("C48").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.643'
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = '  +6.09%  '
$ws.Range("B49").Value = 'VeChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0259'
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = '  +7.38%  '
$ws.Range("E50").Value = '  +1.86%  '
$ws.Range("B51").Value = 'Maker'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D51").Value = '2.171.78'
$ws.Range("E51").Value = '  +7.50%  '
